# Auto-generated edit script: apply numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# across all 8 sheets, per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 868.9231
$ws.Range("I9").Value = 314.14285
$ws.Range("J9").Value = 1516.1666
$ws.Range("K9").Value = 314.14285
$ws.Range("L9").Value = 1516.1666
$ws.Range("M9").Value = -145.14285
$ws.Range("N9").Value = -1854.1666
$ws.Range("H12").Value = 336
$ws.Range("I12").Value = 412.75
$ws.Range("K12").Value = 412.75
$ws.Range("M12").Value = -242.75
$ws.Range("H42").Value = 1463.3334
$ws.Range("I42").Value = 536
$ws.Range("K42").Value = 1608
$ws.Range("M42").Value = -1378
$ws.Range("H51").Value = 4633.1113
$ws.Range("J51").Value = 4633.1113
$ws.Range("L51").Value = 4633.1113
$ws.Range("N51").Value = -5601.1113
$ws.Range("H69").Value = 14582.941
$ws.Range("J69").Value = 14582.941
$ws.Range("L69").Value = 43748.823
$ws.Range("N69").Value = -45496.823
$ws.Range("H72").Value = 14582.941
$ws.Range("J72").Value = 14582.941
$ws.Range("L72").Value = 131246.469
$ws.Range("N72").Value = -139982.469
$ws.Range("H74").Value = 33186.285
$ws.Range("I74").Value = 27279.785
$ws.Range("J74").Value = 44999.285
$ws.Range("K74").Value = 27279.785
$ws.Range("L74").Value = 44999.285
$ws.Range("M74").Value = -26343.785
$ws.Range("N74").Value = -46871.285
$ws.Range("H77").Value = 33186.285
$ws.Range("I77").Value = 27279.785
$ws.Range("J77").Value = 44999.285
$ws.Range("K77").Value = 136398.925
$ws.Range("L77").Value = 224996.425
$ws.Range("M77").Value = -131718.925
$ws.Range("N77").Value = -234356.425
$ws.Range("H86").Value = 7666
$ws.Range("J86").Value = 7666
$ws.Range("L86").Value = 7666
$ws.Range("N86").Value = -9912
$ws.Range("H89").Value = 7666
$ws.Range("J89").Value = 7666
$ws.Range("L89").Value = 38330
$ws.Range("N89").Value = -49562
$ws.Range("H92").Value = 610.73334
$ws.Range("I92").Value = 547.25
$ws.Range("J92").Value = 864.6667
$ws.Range("K92").Value = 547.25
$ws.Range("L92").Value = 864.6667
$ws.Range("M92").Value = 700.75
$ws.Range("N92").Value = -3360.6667
$ws.Range("H93").Value = 21999.5
$ws.Range("J93").Value = 21999.5
$ws.Range("L93").Value = 21999.5
$ws.Range("N93").Value = -26991.5
$ws.Range("H95").Value = 28000
$ws.Range("J95").Value = 28000
$ws.Range("L95").Value = 28000
$ws.Range("N95").Value = -33492
$ws.Range("H98").Value = 1833.5807
$ws.Range("I98").Value = 743.7083
$ws.Range("K98").Value = 743.7083
$ws.Range("M98").Value = 754.2917
$ws.Range("H111").Value = 632.9167
$ws.Range("I111").Value = 624.7
$ws.Range("J111").Value = 674
$ws.Range("K111").Value = 1874.1
$ws.Range("L111").Value = 2022
$ws.Range("M111").Value = 1192.9
$ws.Range("N111").Value = -8156
$ws.Range("H114").Value = 90000
$ws.Range("I114").Value = 90000
$ws.Range("K114").Value = 90000
$ws.Range("M114").Value = -85661
$ws.Range("H116").Value = 5831.6665
$ws.Range("J116").Value = 6747.75
$ws.Range("L116").Value = 6747.75
$ws.Range("N116").Value = -13631.75
$ws.Range("H122").Value = 1833.5807
$ws.Range("I122").Value = 743.7083
$ws.Range("K122").Value = 2231.1249
$ws.Range("M122").Value = 218.8751000000002
$ws.Range("H127").Value = 3349.4666
$ws.Range("I127").Value = 2608
$ws.Range("J127").Value = 4461.6665
$ws.Range("K127").Value = 7824
$ws.Range("L127").Value = 13384.9995
$ws.Range("M127").Value = -2864
$ws.Range("N127").Value = -23304.9995
$ws.Range("H135").Value = 1384.7222
$ws.Range("I135").Value = 1261.8
$ws.Range("J135").Value = 1999.3334
$ws.Range("K135").Value = 11356.2
$ws.Range("L135").Value = 17994.0006
$ws.Range("M135").Value = -8821.199999999999
$ws.Range("N135").Value = -23064.0006
$ws.Range("H137").Value = 6277.8335
$ws.Range("J137").Value = 8789.5
$ws.Range("L137").Value = 26368.5
$ws.Range("N137").Value = -31468.5
$ws.Range("H138").Value = 5123.102
$ws.Range("I138").Value = 4938.6
$ws.Range("J138").Value = 5161.0137
$ws.Range("K138").Value = 14815.8
$ws.Range("L138").Value = 15483.0411
$ws.Range("M138").Value = -9675.800000000001
$ws.Range("N138").Value = -25763.0411
$ws.Range("H141").Value = 2523.0625
$ws.Range("I141").Value = 2598.6
$ws.Range("J141").Value = 1390
$ws.Range("K141").Value = 7795.799999999999
$ws.Range("L141").Value = 4170
$ws.Range("M141").Value = -2615.799999999999
$ws.Range("N141").Value = -14530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3147.6667
$ws.Range("I2").Value = 2779.2
$ws.Range("K2").Value = 2779.2
$ws.Range("M2").Value = -2666.2
$ws.Range("H10").Value = 3833
$ws.Range("I10").Value = 749.5
$ws.Range("K10").Value = 749.5
$ws.Range("M10").Value = -579.5
$ws.Range("H32").Value = 13132.757
$ws.Range("I32").Value = 8490.968000000001
$ws.Range("K32").Value = 8490.968000000001
$ws.Range("M32").Value = -8203.968000000001
$ws.Range("H45").Value = 2073.375
$ws.Range("I45").Value = 1778.8334
$ws.Range("K45").Value = 1778.8334
$ws.Range("M45").Value = -1401.8334
$ws.Range("H55").Value = 103749.75
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 103749.75
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 103749.75
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -104379.75
$ws.Range("H63").Value = 1995.2858
$ws.Range("I63").Value = 1995.8334
$ws.Range("K63").Value = 1995.8334
$ws.Range("M63").Value = -1309.8334
$ws.Range("H66").Value = 1995.2858
$ws.Range("I66").Value = 1995.8334
$ws.Range("K66").Value = 9979.166999999999
$ws.Range("M66").Value = -6547.166999999999
$ws.Range("H69").Value = 394999.66
$ws.Range("J69").Value = 394999.66
$ws.Range("L69").Value = 394999.66
$ws.Range("N69").Value = -396497.66
$ws.Range("H72").Value = 394999.66
$ws.Range("J72").Value = 394999.66
$ws.Range("L72").Value = 1184998.98
$ws.Range("N72").Value = -1192486.98
$ws.Range("H88").Value = 1320.7
$ws.Range("J88").Value = 1716.1666
$ws.Range("L88").Value = 1716.1666
$ws.Range("N88").Value = -2528.1666
$ws.Range("H91").Value = 1320.7
$ws.Range("J91").Value = 1716.1666
$ws.Range("L91").Value = 1716.1666
$ws.Range("N91").Value = -4524.1666
$ws.Range("H110").Value = 6002.75
$ws.Range("I110").Value = 7333
$ws.Range("K110").Value = 7333
$ws.Range("M110").Value = -5288
$ws.Range("H116").Value = 3147.6667
$ws.Range("I116").Value = 2779.2
$ws.Range("K116").Value = 2779.2
$ws.Range("M116").Value = -485.1999999999998
$ws.Range("H132").Value = 3098.2856
$ws.Range("I132").Value = 2903.4211
$ws.Range("J132").Value = 4949.5
$ws.Range("K132").Value = 8710.263300000001
$ws.Range("L132").Value = 14848.5
$ws.Range("M132").Value = -6180.263300000001
$ws.Range("N132").Value = -19908.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3147.6667
$ws.Range("I3").Value = 2779.2
$ws.Range("K3").Value = 2779.2
$ws.Range("M3").Value = -2665.2
$ws.Range("H14").Value = 8
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = 8
$ws.Range("M14").Value = 164
$ws.Range("H20").Value = 2517.9375
$ws.Range("I20").Value = 2559.2666
$ws.Range("K20").Value = 2559.2666
$ws.Range("M20").Value = -2312.2666
$ws.Range("H35").Value = 38524.8
$ws.Range("I35").Value = 38524.8
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 38524.8
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -38214.8
$ws.Range("N35").ClearContents()
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21372
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66864
$ws.Range("H86").Value = 5700
$ws.Range("I86").Value = 3275.75
$ws.Range("K86").Value = 3275.75
$ws.Range("M86").Value = -2152.75
$ws.Range("H89").Value = 5700
$ws.Range("I89").Value = 3275.75
$ws.Range("K89").Value = 16378.75
$ws.Range("M89").Value = -10762.75
$ws.Range("H94").Value = 1456.1428
$ws.Range("I94").Value = 838.6
$ws.Range("K94").Value = 838.6
$ws.Range("M94").Value = -387.6
$ws.Range("H97").Value = 5571.2
$ws.Range("I97").Value = 5571.2
$ws.Range("K97").Value = 5571.2
$ws.Range("M97").Value = -4580.2
$ws.Range("H105").Value = 1789.5714
$ws.Range("I105").Value = 1395.0625
$ws.Range("J105").Value = 3052
$ws.Range("K105").Value = 1395.0625
$ws.Range("L105").Value = 3052
$ws.Range("M105").Value = 351.9375
$ws.Range("N105").Value = -6546
$ws.Range("H107").Value = 2181.2222
$ws.Range("I107").Value = 1945.4286
$ws.Range("K107").Value = 1945.4286
$ws.Range("M107").Value = -25.42859999999996
$ws.Range("H134").Value = 2901.5
$ws.Range("I134").Value = 1857.2858
$ws.Range("K134").Value = 5571.857400000001
$ws.Range("M134").Value = -3036.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4198.4287
$ws.Range("I16").Value = 4179.2
$ws.Range("K16").Value = 4179.2
$ws.Range("M16").Value = -3892.2
$ws.Range("H22").Value = 2740.2856
$ws.Range("I22").Value = 319.16666
$ws.Range("K22").Value = 319.16666
$ws.Range("M22").Value = 30.83334000000002
$ws.Range("H31").Value = 7123.12
$ws.Range("I31").Value = 7149.136
$ws.Range("J31").Value = 6932.3335
$ws.Range("K31").Value = 7149.136
$ws.Range("L31").Value = 6932.3335
$ws.Range("M31").Value = -6854.136
$ws.Range("N31").Value = -7522.3335
$ws.Range("H34").Value = 7123.12
$ws.Range("I34").Value = 7149.136
$ws.Range("J34").Value = 6932.3335
$ws.Range("K34").Value = 7149.136
$ws.Range("L34").Value = 6932.3335
$ws.Range("M34").Value = -6947.136
$ws.Range("N34").Value = -7336.3335
$ws.Range("H99").Value = 39241.9
$ws.Range("I99").Value = 52837
$ws.Range("K99").Value = 52837
$ws.Range("M99").Value = -51339
$ws.Range("H107").Value = 517.5789
$ws.Range("J107").Value = 1032.1666
$ws.Range("L107").Value = 1032.1666
$ws.Range("N107").Value = -4872.1666
$ws.Range("H113").Value = 4198.4287
$ws.Range("I113").Value = 4179.2
$ws.Range("K113").Value = 4179.2
$ws.Range("M113").Value = -2009.2
$ws.Range("H122").Value = 1929.4
$ws.Range("I122").Value = 1929.4
$ws.Range("K122").Value = 5788.200000000001
$ws.Range("M122").Value = -3338.200000000001
$ws.Range("H126").Value = 39241.9
$ws.Range("I126").Value = 52837
$ws.Range("K126").Value = 158511
$ws.Range("M126").Value = -156041
$ws.Range("H132").Value = 339
$ws.Range("I132").Value = 339
$ws.Range("K132").Value = 1017
$ws.Range("M132").Value = 1513
$ws.Range("H134").Value = 3615
$ws.Range("I134").Value = 3615
$ws.Range("K134").Value = 10845
$ws.Range("M134").Value = -8310
$ws.Range("H141").Value = 366425.22
$ws.Range("J141").Value = 366425.22
$ws.Range("L141").Value = 366425.22
$ws.Range("N141").Value = -376785.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100530.45
$ws.Range("I2").Value = 64733
$ws.Range("J2").Value = 222241.8
$ws.Range("K2").Value = 388398
$ws.Range("L2").Value = 1333450.8
$ws.Range("M2").Value = -388285
$ws.Range("N2").Value = -1333676.8
$ws.Range("H34").Value = 3083.3333
$ws.Range("J34").Value = 3083.3333
$ws.Range("L34").Value = 9249.999899999999
$ws.Range("N34").Value = -9417.999899999999
$ws.Range("H50").Value = 933.3
$ws.Range("I50").Value = 1505.5
$ws.Range("J50").Value = 551.8333
$ws.Range("K50").Value = 4516.5
$ws.Range("L50").Value = 1655.4999
$ws.Range("M50").Value = -4035.5
$ws.Range("N50").Value = -2617.4999
$ws.Range("H53").Value = 933.3
$ws.Range("I53").Value = 1505.5
$ws.Range("J53").Value = 551.8333
$ws.Range("K53").Value = 4516.5
$ws.Range("L53").Value = 1655.4999
$ws.Range("M53").Value = -4035.5
$ws.Range("N53").Value = -2617.4999
$ws.Range("H92").Value = 979.2
$ws.Range("I92").Value = 999
$ws.Range("J92").Value = 949.5
$ws.Range("K92").Value = 2997
$ws.Range("L92").Value = 2848.5
$ws.Range("M92").Value = -1749
$ws.Range("N92").Value = -5344.5
$ws.Range("H113").Value = 2362.9565
$ws.Range("J113").Value = 2352.95
$ws.Range("L113").Value = 7058.849999999999
$ws.Range("N113").Value = -11398.85
$ws.Range("H129").Value = 2055.5715
$ws.Range("I129").Value = 977.8
$ws.Range("J129").Value = 4750
$ws.Range("K129").Value = 2933.4
$ws.Range("L129").Value = 14250
$ws.Range("M129").Value = 2066.6
$ws.Range("N129").Value = -24250
$ws.Range("H139").Value = 55558944
$ws.Range("J139").Value = 3415.3333
$ws.Range("L139").Value = 10245.9999
$ws.Range("N139").Value = -20525.9999
$ws.Range("H140").Value = 11408.909
$ws.Range("J140").Value = 14187.5
$ws.Range("L140").Value = 42562.5
$ws.Range("N140").Value = -52922.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20368
$ws.Range("H70").Value = 5789.231
$ws.Range("J70").Value = 5897.5454
$ws.Range("L70").Value = 5897.5454
$ws.Range("N70").Value = -6437.5454
$ws.Range("H73").Value = 5789.231
$ws.Range("J73").Value = 5897.5454
$ws.Range("L73").Value = 5897.5454
$ws.Range("N73").Value = -7769.5454
$ws.Range("H102").Value = 1517
$ws.Range("I102").Value = 1479.8
$ws.Range("K102").Value = 1479.8
$ws.Range("M102").Value = 142.2
$ws.Range("H104").Value = 39999
$ws.Range("J104").Value = 39999
$ws.Range("L104").Value = 39999
$ws.Range("N104").Value = -46987
$ws.Range("H107").Value = 536.7143
$ws.Range("I107").Value = 467.83334
$ws.Range("K107").Value = 467.83334
$ws.Range("M107").Value = 1452.16666
$ws.Range("H132").Value = 3066.8
$ws.Range("I132").Value = 2963.111
$ws.Range("K132").Value = 8889.332999999999
$ws.Range("M132").Value = -6359.332999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 299.5
$ws.Range("I19").Value = 299.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 299.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -129.5
$ws.Range("N19").ClearContents()
$ws.Range("H46").Value = 3041.0557
$ws.Range("I46").Value = 2694.2727
$ws.Range("J46").Value = 3586
$ws.Range("K46").Value = 2694.2727
$ws.Range("L46").Value = 3586
$ws.Range("M46").Value = -2506.2727
$ws.Range("N46").Value = -3962
$ws.Range("H61").Value = 19341.166
$ws.Range("I61").Value = 21566.111
$ws.Range("K61").Value = 21566.111
$ws.Range("M61").Value = -21364.111
$ws.Range("H93").Value = 1194.25
$ws.Range("I93").Value = 700
$ws.Range("J93").Value = 1359
$ws.Range("K93").Value = 700
$ws.Range("L93").Value = 1359
$ws.Range("M93").Value = 548
$ws.Range("N93").Value = -3855
$ws.Range("H113").Value = 19341.166
$ws.Range("I113").Value = 21566.111
$ws.Range("K113").Value = 21566.111
$ws.Range("M113").Value = -19396.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 958.3333
$ws.Range("I17").Value = 1350
$ws.Range("J17").Value = 175
$ws.Range("K17").Value = 1350
$ws.Range("L17").Value = 175
$ws.Range("M17").Value = -1178
$ws.Range("N17").Value = -519
$ws.Range("H50").Value = 32084
$ws.Range("J50").Value = 32084
$ws.Range("L50").Value = 32084
$ws.Range("N50").Value = -33346
$ws.Range("H51").Value = 40000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9490
$ws.Range("H107").Value = 722.8
$ws.Range("I107").Value = 636.55554
$ws.Range("K107").Value = 1909.66662
$ws.Range("M107").Value = 10.33338000000003
$ws.Range("H113").Value = 838.4
$ws.Range("I113").Value = 498
$ws.Range("K113").Value = 1494
$ws.Range("M113").Value = 676
$ws.Range("H122").Value = 3810
$ws.Range("J122").Value = 5227.5
$ws.Range("L122").Value = 15682.5
$ws.Range("N122").Value = -20582.5
$ws.Range("H136").Value = 3040.6843
$ws.Range("I136").Value = 2962.7222
$ws.Range("K136").Value = 8888.1666
$ws.Range("M136").Value = -6338.1666

